$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the numeric-looking Price (column D) updates to be stored as plain text
# (matching the original inline-string cells) instead of being auto-converted
# to numbers by Excels input parser. We flip column D to Text format, write the
# values, then clear the formatting again so the cells end up with no explicit
# style (same as the source file).
$numericPriceCells = @("D5","D6","D7","D8","D9","D10","D11","D12","D14","D15","D16","D18","D20","D22","D24","D25","D27","D28","D29","D30","D31","D32","D33","D35","D36","D37","D38","D39","D41","D42","D43","D45","D46","D47","D48","D49","D50","D51")
foreach ($addr in $numericPriceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# --- Column D (Price) updates ---
$ws.Range("D5").Value = "330.11"
$ws.Range("D6").Value = "1.003"
$ws.Range("D7").Value = "0.5212"
$ws.Range("D8").Value = "0.4409"
$ws.Range("D9").Value = "54.07"
$ws.Range("D10").Value = "0.08923"
$ws.Range("D11").Value = "1.151"
$ws.Range("D12").Value = "24.22"
$ws.Range("D14").Value = "6.681"
$ws.Range("D15").Value = "7.680"
$ws.Range("D16").Value = "95.83"
$ws.Range("D18").Value = "0.00001121"
$ws.Range("D20").Value = "19.13"
$ws.Range("D22").Value = "6.244"
$ws.Range("D24").Value = "12.28"
$ws.Range("D25").Value = "2.309"
$ws.Range("D27").Value = "22.19"
$ws.Range("D28").Value = "2.551"
$ws.Range("D29").Value = "163.55"
$ws.Range("D30").Value = "131.51"
$ws.Range("D31").Value = "1.186"
$ws.Range("D32").Value = "0.1067"
$ws.Range("D33").Value = "1.648"
$ws.Range("D35").Value = "3.903"
$ws.Range("D36").Value = "10.06"
$ws.Range("D37").Value = "0.02552"
$ws.Range("D38").Value = "0.06802"
$ws.Range("D39").Value = "5.464"
$ws.Range("D41").Value = "0.2248"
$ws.Range("D42").Value = "0.6854"
$ws.Range("D43").Value = "1.250"
$ws.Range("D45").Value = "13.91"
$ws.Range("D46").Value = "0.6311"
$ws.Range("D47").Value = "2.191"
$ws.Range("D48").Value = "3.623"
$ws.Range("D49").Value = "1.236"
$ws.Range("D50").Value = "1.241"
$ws.Range("D51").Value = "81.47"

$ws.Range("D2").Value = "30.481.22"
$ws.Range("D3").Value = "2.092.45"
$ws.Range("D13").Value = "2.098.14"
$ws.Range("D23").Value = "30.512.68"
$ws.Range("D26").Value = "2.336.37"

# Clear the temporary Text formatting so these cells go back to having no
# explicit style, matching the source workbook.
foreach ($addr in $numericPriceCells) {
    $ws.Range($addr).ClearFormats()
}

# --- Column E (Volume/1h) updates ---
$ws.Range("E2").Value = "  -1.33%  "
$ws.Range("E3").Value = "  -1.36%  "
$ws.Range("E4").Value = "  +0.29%  "
$ws.Range("E5").Value = "  -1.55%  "
$ws.Range("E6").Value = "  +0.31%  "
$ws.Range("E7").Value = "  -2.98%  "
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("E9").Value = "  +16.39%  "
$ws.Range("E10").Value = "  -1.73%  "
$ws.Range("E11").Value = "  -2.88%  "
$ws.Range("E12").Value = "  -4.89%  "
$ws.Range("E13").Value = "  -1.24%  "
$ws.Range("E14").Value = "  -1.69%  "
$ws.Range("E15").Value = "  -2.24%  "
$ws.Range("E16").Value = "  -2.56%  "
$ws.Range("E17").Value = "  +0.35%  "
$ws.Range("E18").Value = "  -1.86%  "
$ws.Range("E19").Value = "  -0.76%  "
$ws.Range("E20").Value = "  -0.85%  "
$ws.Range("E21").Value = "  +0.24%  "
$ws.Range("E22").Value = "  -2.60%  "
$ws.Range("E23").Value = "  -1.49%  "
$ws.Range("E24").Value = "  +0.66%  "
$ws.Range("E26").Value = "  -1.55%  "
$ws.Range("E27").Value = "  -3.46%  "
$ws.Range("E28").Value = "  -2.13%  "
$ws.Range("E29").Value = "  -0.26%  "
$ws.Range("E30").Value = "  -1.97%  "
$ws.Range("E31").Value = "  +0.65%  "
$ws.Range("E33").Value = "  +5.90%  "
$ws.Range("E34").Value = "  -2.42%  "
$ws.Range("E35").Value = "  -2.51%  "
$ws.Range("E36").Value = "  +3.92%  "
$ws.Range("E37").Value = "  -2.85%  "
$ws.Range("E38").Value = "  +0.38%  "
$ws.Range("E39").Value = "  -2.68%  "
$ws.Range("E40").Value = "  -5.96%  "
$ws.Range("E41").Value = "  -1.54%  "
$ws.Range("E42").Value = "  -0.46%  "
$ws.Range("E43").Value = "  -0.96%  "
$ws.Range("E44").Value = "  +0.25%  "
$ws.Range("E45").Value = "  -2.03%  "
$ws.Range("E46").Value = "  -2.40%  "
$ws.Range("E47").Value = "  -3.53%  "
$ws.Range("E48").Value = "  -1.81%  "
$ws.Range("E49").Value = "  +5.51%  "
$ws.Range("E50").Value = "  -3.94%  "
$ws.Range("E51").Value = "  -2.57%  "
